$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing rows (45-50) that no longer exist in the updated TI domain table
$ws.Range("A45:H50").EntireRow.Delete() | Out-Null

# Update IETESTCD (C), IETEST (D), and IECAT (E) for each remaining data row (2-44)
$ws.Range("C2").Value = 'INCL1'
$ws.Range("D2").Value = 'Inclusion Criteria All patients and subjects will be willing to commit to training and data collection that includes video recording, which may be used for educational and (As per the protocol)'
$ws.Range("E2").Value = 'INCLUSION'
$ws.Range("C3").Value = 'INCL2'
$ws.Range("D3").Value = 'Able to give informed consent and willing to comply with the requirements of the study protocol'
$ws.Range("E3").Value = 'INCLUSION'
$ws.Range("C4").Value = 'INCL3'
$ws.Range("D4").Value = 'Age ≥ 18 years'
$ws.Range("E4").Value = 'INCLUSION'
$ws.Range("C5").Value = 'INCL4'
$ws.Range("D5").Value = 'Be a current resident within the United States'
$ws.Range("E5").Value = 'INCLUSION'
$ws.Range("C6").Value = 'INCL5'
$ws.Range("D6").Value = 'Have had RA for ≥ 6 months, as diagnosed by a qualiﬁed rheumatologist, according to the revised 1987 American College of Rheumatology (ACR) criteria (Arnett et al. 1988)'
$ws.Range("E6").Value = 'INCLUSION'
$ws.Range("C7").Value = 'INCL6'
$ws.Range("D7").Value = 'Patients must be deemed suitable candidates to use an AI at home, in the investigator’s judgment, either by self-administration or from a CG or from a HCP'
$ws.Range("E7").Value = 'INCLUSION'
$ws.Range("C8").Value = 'INCL7'
$ws.Range("D8").Value = 'Have been receiving 162 mg TCZ SC q2w or qw using the commercially Molecule Name and Protocol Name – Organization Name 6 / Protocol XX12346, Version 2.1 available PFS-NSD for at (As per the protocol)'
$ws.Range("E8").Value = 'INCLUSION'
$ws.Range("C9").Value = 'INCL8'
$ws.Range("D9").Value = 'Most recent laboratory results performed in accordance with the current Actemra U.S. Prescribing Information (USPI) do not warrant dose adjustment or discontinuation of therapy (As per the protocol)'
$ws.Range("E9").Value = 'INCLUSION'
$ws.Range("C10").Value = 'INCL9'
$ws.Range("D10").Value = 'At least 2 self-injecting patients will be left-hand dominant'
$ws.Range("E10").Value = 'INCLUSION'
$ws.Range("C11").Value = 'INCL10'
$ws.Range("D11").Value = 'To continue using contraception as discussed with the patient’s rheumatologist at the time of prescription of TCZ SC PFS-NSD 4.1.1.2 Caregivers CGs may be already acquainted (As per the protocol)'
$ws.Range("E11").Value = 'INCLUSION'
$ws.Range("C12").Value = 'INCL11'
$ws.Range("D12").Value = 'Able to give informed consent and willing to comply with the requirements of the study protocol'
$ws.Range("E12").Value = 'INCLUSION'
$ws.Range("C13").Value = 'INCL12'
$ws.Range("D13").Value = 'Age ≥ 18 years'
$ws.Range("E13").Value = 'INCLUSION'
$ws.Range("C14").Value = 'INCL13'
$ws.Range("D14").Value = 'Not professionally qualiﬁed to give an injection (e.g., a patient’s spouse, relative).'
$ws.Range("E14").Value = 'INCLUSION'
$ws.Range("C15").Value = 'INCL14'
$ws.Range("D15").Value = 'Able (after training) and willing to inject a patient at each visit'
$ws.Range("E15").Value = 'INCLUSION'
$ws.Range("C16").Value = 'INCL15'
$ws.Range("D16").Value = 'Must be current resident within the United States 4.1.1.3 Healthcare Professionals An HCP may be already acquainted with (and already supporting) a specific patient enrolled in (As per the protocol)'
$ws.Range("E16").Value = 'INCLUSION'
$ws.Range("C17").Value = 'INCL16'
$ws.Range("D17").Value = 'Able to give informed consent and willing to comply with the requirements of the study protocol'
$ws.Range("E17").Value = 'INCLUSION'
$ws.Range("C18").Value = 'INCL17'
$ws.Range("D18").Value = 'Age ≥ 18 years'
$ws.Range("E18").Value = 'INCLUSION'
$ws.Range("C19").Value = 'INCL18'
$ws.Range("D19").Value = 'Must be current resident of the United States'
$ws.Range("E19").Value = 'INCLUSION'
$ws.Range("C20").Value = 'INCL19'
$ws.Range("D20").Value = 'Professionally qualiﬁed to give an injection and willing to inject a patient and comply with the study protocol'
$ws.Range("E20").Value = 'INCLUSION'
$ws.Range("C21").Value = 'EXCL1'
$ws.Range("D21").Value = 'Exclusion Criteria 4.1.2.1 Patients Patients who meet any of the following criteria will be excluded from study entry:'
$ws.Range("E21").Value = 'EXCLUSION'
$ws.Range("C22").Value = 'EXCL2'
$ws.Range("D22").Value = 'Any serious medical condition or abnormality in clinical laboratory tests that, in the investigator’s judgment, precludes the patient’s safe participation Molecule Name and (As per the protocol)'
$ws.Range("E22").Value = 'EXCLUSION'
$ws.Range("C23").Value = 'EXCL3'
$ws.Range("D23").Value = 'Patients with functional RA status class IV (according to the 1991 ACR revised criteria for the classiﬁcation of global functional status in RA [Hochberg et al. 1992]), as (As per the protocol)'
$ws.Range("E23").Value = 'EXCLUSION'
$ws.Range("C24").Value = 'EXCL4'
$ws.Range("D24").Value = 'Neuropathies or other conditions that might interfere with pain evaluation'
$ws.Range("E24").Value = 'EXCLUSION'
$ws.Range("C25").Value = 'EXCL5'
$ws.Range("D25").Value = 'Current participation in any interventional clinical trial'
$ws.Range("E25").Value = 'EXCLUSION'
$ws.Range("C26").Value = 'EXCL6'
$ws.Range("D26").Value = 'Patients who self-report to be pregnant or nursing (breastfeeding)'
$ws.Range("E26").Value = 'EXCLUSION'
$ws.Range("C27").Value = 'EXCL7'
$ws.Range("D27").Value = 'Patient or anyone in his/her immediate household is employed in the pharmaceutical industry'
$ws.Range("E27").Value = 'EXCLUSION'
$ws.Range("C28").Value = 'EXCL8'
$ws.Range("D28").Value = 'Patient employed by Roche, Genentech, Battelle, or a contract research organization (CRO) involved in this study (WA29917)'
$ws.Range("E28").Value = 'EXCLUSION'
$ws.Range("C29").Value = 'EXCL9'
$ws.Range("D29").Value = 'Participation in any previous Actemra research study that involved an AI.'
$ws.Range("E29").Value = 'EXCLUSION'
$ws.Range("C30").Value = 'EXCL10'
$ws.Range("D30").Value = 'Prior use of the AI-1000 G1 or AI-1000 G2 in any HF study.'
$ws.Range("E30").Value = 'EXCLUSION'
$ws.Range("C31").Value = 'EXCL11'
$ws.Range("D31").Value = 'ANC < 1.0 × 109/L (1000/mm3) at last (as per the USPI) laboratory assessment.'
$ws.Range("E31").Value = 'EXCLUSION'
$ws.Range("C32").Value = 'EXCL12'
$ws.Range("D32").Value = 'Platelet count < 100 × 109/L (100,000/mm3) at last laboratory assessment.'
$ws.Range("E32").Value = 'EXCLUSION'
$ws.Range("C33").Value = 'EXCL13'
$ws.Range("D33").Value = 'ALT or AST > upper limit of normal [ULN] at last laboratory assessment. 4.1.2.2 Caregivers CGs who meet any of the following criteria will be excluded from study entry:'
$ws.Range("E33").Value = 'EXCLUSION'
$ws.Range("C34").Value = 'EXCL14'
$ws.Range("D34").Value = 'Current participation in any interventional clinical trial'
$ws.Range("E34").Value = 'EXCLUSION'
$ws.Range("C35").Value = 'EXCL15'
$ws.Range("D35").Value = 'Subject or anyone in his/her immediate household is employed in the pharmaceutical industry'
$ws.Range("E35").Value = 'EXCLUSION'
$ws.Range("C36").Value = 'EXCL16'
$ws.Range("D36").Value = 'Subject employed by Roche, Genentech, Battelle, or a CRO involved in this study (WA29917) 4.1.2.3 Healthcare Professionals HCPs who meet any of the following criteria will be (As per the protocol)'
$ws.Range("E36").Value = 'EXCLUSION'
$ws.Range("C37").Value = 'EXCL17'
$ws.Range("D37").Value = 'Current participation in any interventional clinical trial as a patient'
$ws.Range("E37").Value = 'EXCLUSION'
$ws.Range("C38").Value = 'EXCL18'
$ws.Range("D38").Value = 'Participation in the conduct or oversight of this study (WA29917)'
$ws.Range("E38").Value = 'EXCLUSION'
$ws.Range("C39").Value = 'EXCL19'
$ws.Range("D39").Value = 'Subject or anyone in his/her immediate household is employed in the pharmaceutical industry.'
$ws.Range("E39").Value = 'EXCLUSION'
$ws.Range("C40").Value = 'EXCL20'
$ws.Range("D40").Value = 'Subject employed by Roche, Genentech, Battelle or a CRO involved in this study (WA29917) Molecule Name and Protocol Name – Organization Name 8 / Protocol XX12346, Version 2.1'
$ws.Range("E40").Value = 'EXCLUSION'
$ws.Range("C41").Value = 'EXCL21'
$ws.Range("D41").Value = 'Not professionally qualiﬁed to give injections 4.2 METHOD OF TREATMENT ASSIGNMENT AND BLINDING Patients will take part in the study in one of three groups:'
$ws.Range("E41").Value = 'EXCLUSION'
$ws.Range("C42").Value = 'EXCL22'
$ws.Range("D42").Value = 'Self-injecting patients'
$ws.Range("E42").Value = 'EXCLUSION'
$ws.Range("C43").Value = 'EXCL23'
$ws.Range("D43").Value = 'Patients who receive injections from a CG'
$ws.Range("E43").Value = 'EXCLUSION'
$ws.Range("C44").Value = 'EXCL24'
$ws.Range("D44").Value = 'Patients who receive injections from an HCP A minimum of 15 patients and/or subjects are required per group, and patients may be assigned so that all groups have the required (As per the protocol)'
$ws.Range("E44").Value = 'EXCLUSION'
